# Update posts.xlsx after post
# The post corresponding to row 355 ("「お金では買えないよ」...") was removed.
# Deleting the entire row shifts every subsequent row up by one, which matches
# the target diff (old row 356 -> new row 355, ..., old row 380 -> new row 379)
# and shrinks the used range from A1:C380 to A1:C379.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(355).Delete()
